$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column D with header "A" (same bold/centered/bordered style as B1/C1) ---
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("D1").Value = "A"

# --- Update row labels ---
$ws.Range("A4").Value = "A Lag"

# --- Update the regression coefficients (columns B and C) ---
$ws.Range("B2").Value = "0.289***"
$ws.Range("C2").Value = "10.423***"

$ws.Range("B3").Value = "-0.024***"
$ws.Range("C3").Value = "-0.249***"

$ws.Range("B4").Value = "0.347***"
$ws.Range("C4").Value = "4.355***"

# --- Populate new column D data rows ---
$ws.Range("D2").Value = "-0.17**"
$ws.Range("D3").Value = "0.025***"
$ws.Range("D4").Value = "-0.648***"

# --- Remove old row 5 (r2_adj) entirely ---
$ws.Rows.Item(5).Delete()
